$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.218.54"
$ws.Range("E2").Value = "'  +0.07%  "

$ws.Range("D3").Value = "'1.905.28"
$ws.Range("E3").Value = "'  +0.02%  "

$ws.Range("E4").Value = "'  +0.30%  "

$ws.Range("D5").Value = "'307.64"
$ws.Range("E5").Value = "'  +0.30%  "

$ws.Range("D6").Value = "'1.003"

$ws.Range("D7").Value = "'0.5261"
$ws.Range("E7").Value = "'  +0.13%  "

$ws.Range("D8").Value = "'0.3811"
$ws.Range("E8").Value = "'  +0.98%  "

$ws.Range("E9").Value = "'  +0.70%  "

$ws.Range("D10").Value = "'21.67"
$ws.Range("E10").Value = "'  +2.42%  "

$ws.Range("E11").Value = "'  +0.51%  "

$ws.Range("D12").Value = "'0.08047"
$ws.Range("E12").Value = "'  -4.19%  "

$ws.Range("D13").Value = "'95.87"
$ws.Range("E13").Value = "'  +0.87%  "

$ws.Range("D14").Value = "'5.366"
$ws.Range("E14").Value = "'  +1.64%  "

$ws.Range("D15").Value = "'1.786.47"
$ws.Range("E15").Value = "'  -6.20%  "

$ws.Range("E16").Value = "'  +0.24%  "

$ws.Range("D17").Value = "'0.000008684"
$ws.Range("E17").Value = "'  +0.89%  "

$ws.Range("D18").Value = "'14.74"

$ws.Range("D19").Value = "'1.003"
$ws.Range("E19").Value = "'  +0.26%  "

$ws.Range("D20").Value = "'27.255.02"
$ws.Range("E20").Value = "'  +0.06%  "

$ws.Range("D21").Value = "'5.126"
$ws.Range("E21").Value = "'  +1.15%  "

$ws.Range("D22").Value = "'10.82"
$ws.Range("E22").Value = "'  +1.88%  "

$ws.Range("D23").Value = "'6.474"
$ws.Range("E23").Value = "'  +0.61%  "

$ws.Range("D24").Value = "'2.354"
$ws.Range("E24").Value = "'  +2.95%  "

$ws.Range("D25").Value = "'149.48"
$ws.Range("E25").Value = "'  +1.35%  "

$ws.Range("D26").Value = "'18.27"
$ws.Range("E26").Value = "'  +0.48%  "

$ws.Range("E27").Value = "'  -0.62%  "

$ws.Range("E28").Value = "'  +1.75%  "

$ws.Range("D29").Value = "'4.844"
$ws.Range("E29").Value = "'  +0.45%  "

$ws.Range("D30").Value = "'4.901"
$ws.Range("E30").Value = "'  -0.57%  "

$ws.Range("D31").Value = "'0.09249"
$ws.Range("E31").Value = "'  -0.36%  "

$ws.Range("D32").Value = "'0.05088"
$ws.Range("E32").Value = "'  +0.34%  "

$ws.Range("D33").Value = "'0.8033"
$ws.Range("E33").Value = "'  -0.65%  "

$ws.Range("D34").Value = "'1.229"
$ws.Range("E34").Value = "'  -0.81%  "

$ws.Range("D35").Value = "'2.973"
$ws.Range("E35").Value = "'  +0.49%  "

$ws.Range("D36").Value = "'3.387"
$ws.Range("E36").Value = "'  +0.15%  "

$ws.Range("D37").Value = "'2.679"
$ws.Range("E37").Value = "'  +2.42%  "

$ws.Range("D38").Value = "'0.5732"
$ws.Range("E38").Value = "'  -0.07%  "

$ws.Range("D39").Value = "'0.01991"
$ws.Range("E39").Value = "'  +0.17%  "

$ws.Range("D40").Value = "'1.086"
$ws.Range("E40").Value = "'  +1.04%  "

$ws.Range("D41").Value = "'9.002"
$ws.Range("E41").Value = "'  +0.24%  "

$ws.Range("D42").Value = "'6.599"
$ws.Range("E42").Value = "'  -0.70%  "

$ws.Range("D43").Value = "'116.64"
$ws.Range("E43").Value = "'  -0.73%  "

$ws.Range("E44").Value = "'  +0.39%  "

$ws.Range("D45").Value = "'0.4917"
$ws.Range("E45").Value = "'  +1.27%  "

$ws.Range("B46").Value = "'PaxDollar"
$ws.Range("C46").Value = "'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").Value = "'1.003"
$ws.Range("E46").Value = "'  +0.21%  "

$ws.Range("B47").Value = "'EnergySwap"
$ws.Range("C47").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'10.15"
$ws.Range("E47").Value = "'  +0.25%  "

$ws.Range("D48").Value = "'1.643"
$ws.Range("E48").Value = "'  +1.59%  "

$ws.Range("D49").Value = "'38.63"
$ws.Range("E49").Value = "'  +3.01%  "

$ws.Range("D50").Value = "'64.45"
$ws.Range("E50").Value = "'  +0.99%  "

$ws.Range("D51").Value = "'0.05963"
$ws.Range("E51").Value = "'  +0.36%  "
